# Remove the trailing "Ver no Jupiter ..." / "(c) 2020 ..." footer block
# (and the blank paragraph that precedes it) that followed the
# "LOT2023: Processos Bioquimicos Industriais (Requisito fraco)" line.

$d = $word.ActiveDocument

# Anchor on the "Ver no Jupiter..." paragraph via Find so we don't depend
# on hard-coded paragraph indices.
$findRange = $d.Content
$found = $findRange.Find.Execute("Ver no Jupiter Salvar em pdf Salvar em docx", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $verPara = $findRange.Paragraphs(1)

    # The blank paragraph right before "Ver no Jupiter..." and the
    # copyright paragraph right after it are removed together with it.
    $prevPara = $verPara.Previous()
    $nextPara = $verPara.Next()

    $deleteRange = $d.Range($prevPara.Range.Start, $nextPara.Range.End)
    $deleteRange.Delete()
}
